$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update "Objetivos:" row (row 10) with the full syllabus objectives text
#    (both the B and C "comparison" columns get the same text).
# ---------------------------------------------------------------------------
$objetivosText = "Sistemas de classificação dos aços, estudo da influência dos elementos de liga, características, propriedades e efeito do tratamento térmico das ligas ferrosas.Estudo das ligas de metais não-ferrosos quanto à sua nomenclatura, composição, propriedades, características, tratamentos e aplicações práticas."
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# ---------------------------------------------------------------------------
# 2. Insert a new row at position 13 (pushes old rows 13..23 down to 14..24,
#    carrying their row heights/content along automatically).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The inserted row picked up the column-default style in col A; remove it so
# the row has no A13 cell at all, matching the target layout.
$ws.Range("A13").Clear()

# Populate the new row 13 (B13/C13) with the "Docentes responsaveis" value,
# copying number/cell formatting from an existing B/C data cell first so the
# new cells end up with the correct (wrap-text) styles instead of the
# plain column-default style.
$ws.Range("B9").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$docentesText = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("B13").Value = $docentesText
$ws.Range("C13").Value = $docentesText

# ---------------------------------------------------------------------------
# 3. Replace "Programa resumido:" / "Programa:" values (now rows 14 and 16)
#    with the new short-syllabus text.
# ---------------------------------------------------------------------------
$resumidoText = "- Aços carbono, microligados e inoxidáveis- Ligas de alumínio- Ligas de cobre- Ligas de níquel- Ligas de titânio- Ligas especiais"
$ws.Range("B14").Value = $resumidoText
$ws.Range("C14").Value = $resumidoText
$ws.Range("B16").Value = $resumidoText
$ws.Range("C16").Value = $resumidoText

# ---------------------------------------------------------------------------
# 4. Shift the evaluation-related values (Metodo/Criterio/Norma) up one slot
#    relative to where the row-insert left them, i.e. re-assign the trio of
#    texts to rows 19, 20 and 21.
# ---------------------------------------------------------------------------
$metodoText = "Duas avaliações ao longo do semestre (P1 e P2). A P1 será ministrada por volta da metade do semestre. A P2 abrangerá toda a matéria ministrada no semestre"
$criterioText = "A nota final (NF) será a média aritmética das duas avaliações:NF = (P1 + P2)/2"
$normaText = "Uma única avaliação escrita contemplando toda a matéria ministrada no semestre letivo."

$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# ---------------------------------------------------------------------------
# 5. Fill in the bibliography text on row 22 (B22/C22).
# ---------------------------------------------------------------------------
$biblioText = "1. Bresciani Fº, E. Seleção de materiais metálicos. Editora da UNICAMP,3ª ed., 1991.2. Bresciani Fº, E. Seleção de metais não-ferrosos. Editora da UNICAMP, 1ª ed., 1992.3. Briggs, C.M. Steel Castings Handbook, Steel Founders Society of America, 1970.4. Coutinho, .B. Materiais Metálicos para Engenharia, Fundação Cristiano Ottoni, UFMG, 1992.5. Douglass, D.L. The Metallurgy of Zirconium, International Atomic Energy Agency, 1971.6. Industrial Applications of Titanium and Zirconium, Webster & Young, 1984.7. Donachie Jr., M.J. Titanium: A Technical Guide, ASM International, 1988.8. Jaffe,  R.I. & Promisel, N.G... The Science, Technology and Application of Titanium, Pergamon Press, 1970.9. Cerqueira Leite, R.C. e outros. Nióbio: Uma Conquista Naciuonal, Livraria Duas Cidades, 1988.10. Douglass, D.L. & Kunz, F.W. Columbium Metallurgy, Interscience Publishers, 1961.11. King F. Aluminium and Its Alloys, Ellis Horwood, 1987.12. Hatch, J.E. Aluminium: Properties and Physical Metallurgy1984.13. Dawson, R.J.C. Fusiun Welding and Brazing of Copper and Copper Alloys, Newnes-Butterworths.14. Butts, A. Copper: The Science and Technology of Metals, Its Alloys and Compounds, Reihold Publising, 1954.15. Betterridge, W. Nickel and Its Alloys, McDonald and Evans.16. Hampel, C.A. Rare Metals Handbook, Reinhold Publishing, 1954.17. Yih, S.W. Tungsten: Sources, Metallurgy, Properties and Applications, Plenum Press, 1979.18.Metals Handbook, 9ª ed., vol 1,2 e 3, ASM, 1979.19.Buck, R.M.; Pense, A.W.; Gordon, R.B. Struture and Properties of Engineering Materials, McGraw-Hill, 1977."

$ws.Range("B22").Value = $biblioText
$ws.Range("C22").Value = $biblioText

